$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the existing (currently empty) "Aufgabe 7: " paragraph and
# collapse the found range to right after the label text, i.e. just
# before its paragraph mark.
# ------------------------------------------------------------------
$target = $d.Content
$ok = $target.Find.Execute("Aufgabe 7: ", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find the 'Aufgabe 7: ' paragraph to extend."
}
$target.Collapse(0)

# ------------------------------------------------------------------
# Append the new Aufgabe 7 SQL statement text right after the label.
# A trailing sentinel character ("|") is appended too -- it only
# exists temporarily so that the point where the "_GoBack" bookmark
# must be dropped (see below) is not the literal last character of
# the paragraph: placing a bookmark exactly on a paragraph's
# trailing paragraph-mark position is mishandled by this host and
# resets it to the start of the document. The sentinel is removed
# again immediately after the bookmark has been placed.
# ------------------------------------------------------------------
$target.InsertAfter("SELECT Lieferant, Artikelname, Angebotsdatum FROM liefangebot LEFT JOIN artikel ON Artikel = Artikelnr WHERE DATEDIFF('2017-10-01', Angebotsdatum) > 15|")
$insertEnd = $target.End

# ------------------------------------------------------------------
# The document carries a special "_GoBack" bookmark -- Word's
# last-edit-location marker. It used to sit at the end of the
# "Aufgabe 6" paragraph; re-adding a bookmark with that reserved
# name relocates it (exactly like Word itself does whenever you
# edit elsewhere), landing it right after the text we just typed
# into the Aufgabe 7 paragraph (i.e. just before the sentinel).
# ------------------------------------------------------------------
$bookmarkRange = $d.Range($insertEnd - 1, $insertEnd - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# Remove the temporary sentinel character again.
$sentinelRange = $d.Range($insertEnd - 1, $insertEnd)
$sentinelRange.Delete()
